# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly scraped counts (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (sheet1) ---
$ws1.Range("F2").Value  = 38
$ws1.Range("F7").Value  = 4537
$ws1.Range("F11").Value = 87
$ws1.Range("F14").Value = 177
$ws1.Range("F15").Value = 959
$ws1.Range("F17").Value = 237
$ws1.Range("F22").Value = 3450
$ws1.Range("F23").Value = 5790
$ws1.Range("F29").Value = 3344
$ws1.Range("F34").Value = 517
$ws1.Range("F37").Value = 258
$ws1.Range("F38").Value = 346
$ws1.Range("F40").Value = 1008
$ws1.Range("F41").Value = 897
$ws1.Range("F43").Value = 15
$ws1.Range("F45").Value = 41
$ws1.Range("F46").Value = 465
$ws1.Range("F48").Value = 547

# --- Sheet "全部类型" (sheet4) ---
$ws4.Range("F2").Value  = 38
$ws4.Range("F7").Value  = 4537
$ws4.Range("F12").Value = 87
$ws4.Range("F15").Value = 177
$ws4.Range("F16").Value = 959
$ws4.Range("F18").Value = 237
$ws4.Range("F23").Value = 3450
$ws4.Range("F24").Value = 5790
$ws4.Range("F30").Value = 3344
$ws4.Range("F35").Value = 517
$ws4.Range("F38").Value = 258
$ws4.Range("F39").Value = 346
$ws4.Range("F41").Value = 1008
$ws4.Range("F42").Value = 897
$ws4.Range("F44").Value = 15
$ws4.Range("F46").Value = 41
$ws4.Range("F47").Value = 465
$ws4.Range("F49").Value = 547

$wb.Save()
